# Add a new "Recommended video links (share as homework)" slide as the
# last slide of the deck (position 16), using the Blank layout like every
# other slide in this presentation.

$p = $ppt.ActivePresentation

$newIndex = $p.Slides.Count + 1
$s = $p.Slides.Add($newIndex, 12)   # ppLayoutBlank

# ---------------------------------------------------------------------
# Shape 1: Title textbox
# ---------------------------------------------------------------------
$title = $s.Shapes.AddTextbox(1, 32.4, 11.52, 878.4, 51.84)
$title.Name = "TextBox 1"
$title.Fill.Visible = 0

$ttf = $title.TextFrame
$ttf.TextRange.Text = "Recommended video links (share as homework)"
$ttf.TextRange.ParagraphFormat.Alignment = 1
$ttf.TextRange.Font.Size = 34
$ttf.TextRange.Font.Bold = $true
$ttf.TextRange.Font.Color.RGB = 0x4A2814
$ttf.WordWrap = 0
$ttf.AutoSize = 1
$title.Height = 51.84

# ---------------------------------------------------------------------
# Shape 2: Subtitle textbox
# ---------------------------------------------------------------------
$subtitle = $s.Shapes.AddTextbox(1, 34.56, 59.04, 878.4, 33.12)
$subtitle.Name = "TextBox 2"
$subtitle.Fill.Visible = 0

$stf = $subtitle.TextFrame
$stf.TextRange.Text = "Use these links to support practice after the 45-minute session."
$stf.TextRange.Font.Size = 18
$stf.TextRange.Font.Color.RGB = 0x785644
$stf.WordWrap = 0
$stf.AutoSize = 1
$subtitle.Height = 33.12001   # compensates for float round-trip so EMU lands on 420624

# ---------------------------------------------------------------------
# Shape 3: Left column (links 1-6)
# ---------------------------------------------------------------------
$leftHeadings = @(
    "1. Cursor AI beginner tutorial",
    "2. Cursor MCP setup tutorial",
    "3. Model Context Protocol explained",
    "4. Anthropic MCP tutorial",
    "5. Jira REST API tutorial",
    "6. Figma API tutorial"
)
$leftLinks = @(
    "https://www.youtube.com/results?search_query=Cursor+AI+beginner+tutorial",
    "https://www.youtube.com/results?search_query=Cursor+MCP+setup+tutorial",
    "https://www.youtube.com/results?search_query=Model+Context+Protocol+explained",
    "https://www.youtube.com/results?search_query=Anthropic+MCP+tutorial",
    "https://www.youtube.com/results?search_query=Jira+REST+API+tutorial+developers",
    "https://www.youtube.com/results?search_query=Figma+API+tutorial+for+developers"
)

$left = $s.Shapes.AddTextbox(1, 39.6, 111.6, 439.2, 414.0)
$left.Name = "TextBox 3"
$left.Fill.Visible = 0

$ltf = $left.TextFrame
$ltf.WordWrap = -1

$lines = @()
for ($i = 0; $i -lt $leftHeadings.Count; $i++) {
    $lines += $leftHeadings[$i]
    $lines += $leftLinks[$i]
}
$ltf.TextRange.Text = [string]::Join("`r", $lines)

for ($i = 1; $i -le $lines.Count; $i++) {
    $para = $ltf.TextRange.Paragraphs($i, 1)
    if ($i % 2 -eq 1) {
        $para.ParagraphFormat.SpaceAfter = 0
        $para.Font.Size = 17
        $para.Font.Bold = $true
        $para.Font.Color.RGB = 0x583420
    } else {
        $para.ParagraphFormat.SpaceAfter = 8
        $para.Font.Size = 12
        $para.Font.Color.RGB = 0xAB5818
    }
}

$ltf.AutoSize = 1
$left.Height = 414.0

# ---------------------------------------------------------------------
# Shape 4: Right column (links 7-12)
# ---------------------------------------------------------------------
$rightHeadings = @(
    "7. Bitbucket API tutorial",
    "8. Prompt engineering for developers",
    "9. AI agent workflow tutorials",
    "10. Build MCP server in Python",
    "11. API token security best practices",
    "12. LLM governance for enterprise"
)
$rightLinks = @(
    "https://www.youtube.com/results?search_query=Bitbucket+API+tutorial",
    "https://www.youtube.com/results?search_query=Prompt+engineering+for+software+developers",
    "https://www.youtube.com/results?search_query=AI+agent+workflow+tutorial+developers",
    "https://www.youtube.com/results?search_query=Build+MCP+server+Python+tutorial",
    "https://www.youtube.com/results?search_query=API+token+security+best+practices",
    "https://www.youtube.com/results?search_query=LLM+governance+for+enterprise+teams"
)

$right = $s.Shapes.AddTextbox(1, 486.0, 111.6, 435.6, 414.0)
$right.Name = "TextBox 4"
$right.Fill.Visible = 0

$rtf = $right.TextFrame
$rtf.WordWrap = -1

$rlines = @()
for ($i = 0; $i -lt $rightHeadings.Count; $i++) {
    $rlines += $rightHeadings[$i]
    $rlines += $rightLinks[$i]
}
$rtf.TextRange.Text = [string]::Join("`r", $rlines)

for ($i = 1; $i -le $rlines.Count; $i++) {
    $para = $rtf.TextRange.Paragraphs($i, 1)
    if ($i % 2 -eq 1) {
        $para.ParagraphFormat.SpaceAfter = 0
        $para.Font.Size = 17
        $para.Font.Bold = $true
        $para.Font.Color.RGB = 0x583420
    } else {
        $para.ParagraphFormat.SpaceAfter = 8
        $para.Font.Size = 12
        $para.Font.Color.RGB = 0xAB5818
    }
}

$rtf.AutoSize = 1
$right.Height = 414.0

Write-Output "Added video links slide at index $($s.SlideIndex); total slides = $($p.Slides.Count)"
